# Refactored removing unneeded variables
# - Add a new excluded path row ("path" / "C:\Temp\#TEST IMAGES\folder 1") to the
#   "exclude" sheet.
# - Make "exclude" the active/selected sheet (was "include").
# - Update each sheet's remembered cell selection to A2.
# - Widen column B on "exclude" to fit the new, longer value.

$wb = $excel.ActiveWorkbook

$wsInclude = $wb.Worksheets.Item(1)   # "include"
$wsExclude = $wb.Worksheets.Item(2)   # "exclude"

# Add the new excluded-path data row.
$wsExclude.Range("A2").Value = "path"
$wsExclude.Range("B2").Value = "C:\Temp\#TEST IMAGES\folder 1"

# Resize column B so the longer value fits (matches the original bestFit sizing).
$wsExclude.Columns.Item(2).ColumnWidth = 45.6666666666667

# Move the "include" sheet's remembered selection off of B2 and onto A2.
$wsInclude.Range("A2").Select()

# Switch to the "exclude" sheet and select A2 there too -- it becomes the
# active tab/sheet and its multi-column row selection collapses to A2.
$wsExclude.Activate()
$wsExclude.Range("A2").Select()
